# Update cryptocurrency price and volume(1h) values per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.403.45"
$ws.Range("E2").Value = "'  +0.51%  "

$ws.Range("D3").Value = "'1.848.17"
$ws.Range("E3").Value = "'  -0.46%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.04%  "

$ws.Range("D5").Value = "'233.60"
$ws.Range("E5").Value = "'  +0.51%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  +0.00%  "

$ws.Range("D7").Value = "'0.4671"
$ws.Range("E7").Value = "'  -1.56%  "

$ws.Range("E8").Value = "'  -0.52%  "

$ws.Range("D10").Value = "'1.827.81"
$ws.Range("E10").Value = "'  -1.64%  "

$ws.Range("D11").Value = "'0.07472"
$ws.Range("E11").Value = "'  +0.63%  "

$ws.Range("D12").Value = "'16.31"
$ws.Range("E12").Value = "'  +1.46%  "

$ws.Range("D13").Value = "'4.933"

$ws.Range("D14").Value = "'83.93"
$ws.Range("E14").Value = "'  -1.62%  "

$ws.Range("D15").Value = "'0.6201"
$ws.Range("E15").Value = "'  -2.01%  "

$ws.Range("D16").Value = "'30.346.02"
$ws.Range("E16").Value = "'  +0.39%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "'  -0.04%  "

$ws.Range("D18").Value = "'228.94"
$ws.Range("E18").Value = "'  +1.73%  "

$ws.Range("D19").Value = "'0.000007307"
$ws.Range("E19").Value = "'  +0.00%  "

$ws.Range("D20").Value = "'12.38"
$ws.Range("E20").Value = "'  -3.01%  "

$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "'  +0.14%  "

$ws.Range("D22").Value = "'4.922"
$ws.Range("E22").Value = "'  -3.42%  "

$ws.Range("D23").Value = "'5.871"
$ws.Range("E23").Value = "'  -2.55%  "

$ws.Range("D24").Value = "'166.25"
$ws.Range("E24").Value = "'  -0.67%  "

$ws.Range("D25").Value = "'9.147"
$ws.Range("E25").Value = "'  -0.75%  "

$ws.Range("D26").Value = "'17.83"
$ws.Range("E26").Value = "'  +0.30%  "

$ws.Range("D27").Value = "'1.870"
$ws.Range("E27").Value = "'  +0.78%  "

$ws.Range("D28").Value = "'0.1019"
$ws.Range("E28").Value = "'  -0.43%  "

$ws.Range("D29").Value = "'1.375"

$ws.Range("D30").Value = "'4.083"
$ws.Range("E30").Value = "'  -3.34%  "

$ws.Range("E31").Value = "'  -2.22%  "

$ws.Range("D32").Value = "'0.04873"
$ws.Range("E32").Value = "'  -0.14%  "

$ws.Range("D33").Value = "'1.140"
$ws.Range("E33").Value = "'  -0.55%  "

$ws.Range("D34").Value = "'0.7013"
$ws.Range("E34").Value = "'  -3.44%  "

$ws.Range("D35").Value = "'2.689"
$ws.Range("E35").Value = "'  +0.10%  "

$ws.Range("D36").Value = "'0.01920"
$ws.Range("E36").Value = "'  +0.22%  "

$ws.Range("D37").Value = "'2.660"
$ws.Range("E37").Value = "'  +1.27%  "

$ws.Range("D38").Value = "'0.8650"
$ws.Range("E38").Value = "'  -3.87%  "

$ws.Range("D39").Value = "'105.53"
$ws.Range("E39").Value = "'  -0.23%  "

$ws.Range("D40").Value = "'1.934"
$ws.Range("E40").Value = "'  -2.03%  "

$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "'  +0.61%  "

$ws.Range("D42").Value = "'5.517"
$ws.Range("E42").Value = "'  -0.14%  "

$ws.Range("E43").Value = "'  -1.68%  "

$ws.Range("D44").Value = "'7.073"
$ws.Range("E44").Value = "'  +0.54%  "

$ws.Range("D45").Value = "'61.46"
$ws.Range("E45").Value = "'  +0.23%  "

$ws.Range("D46").Value = "'0.1204"
$ws.Range("E46").Value = "'  +0.01%  "

$ws.Range("D47").Value = "'8.594"
$ws.Range("E47").Value = "'  -2.18%  "

$ws.Range("D48").Value = "'33.38"
$ws.Range("E48").Value = "'  +1.45%  "

$ws.Range("D49").Value = "'0.05542"
$ws.Range("E49").Value = "'  -0.84%  "

$ws.Range("D50").Value = "'1.346"
$ws.Range("E50").Value = "'  -3.92%  "

$ws.Range("D51").Value = "'0.3651"
$ws.Range("E51").Value = "'  -1.46%  "

